$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2's URL (existing hyperlink-styled cell) to the new imgur URL.
$ws.Range("A2").Value = "//i.imgur.com/UxwnWnIb.jpg"

# A3: duplicate of the first imgur URL, normal style (new row).
$ws.Range("A3").Value = "//i.imgur.com/UxwnWnIb.jpg"

# A4: keep its existing "Hipervinculo" style, now carries a value.
$ws.Range("A4").Value = "//i.imgur.com/Or0O3Hob.jpg"

# A5:A11 new rows with remaining imgur URLs.
$ws.Range("A5").Value = "//i.imgur.com/aJTAHiWb.jpg"
$ws.Range("A6").Value = "//i.imgur.com/KpI0LyKb.jpg"
$ws.Range("A7").Value = "//i.imgur.com/fpG0m7Jb.jpg"
$ws.Range("A8").Value = "//i.imgur.com/D9heDsUb.jpg"
$ws.Range("A9").Value = "//i.imgur.com/jZeZz0qb.jpg"
$ws.Range("A10").Value = "//i.imgur.com/iiCBuzvb.jpg"
$ws.Range("A11").Value = "//i.imgur.com/skS4fTJb.jpg"

# Match the final selection shown in the diff.
$ws.Range("D7").Select()
